$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target data for columns B (geneSymbol), C (geneName), D (geneConfidence),
# F (time_taken) for rows 2..108 (in row order). This reflects the insertion of a
# new "COX20" gene row (in the position previously held by COX6A1) plus the new
# time_taken column F with a per-row timestamp.
$data = @(
    @('AARS', 'alanyl-tRNA synthetase', '3', '2021-10-05 10:51:00.539831'),
    @('ATL1', 'atlastin GTPase 1', '3', '2021-10-05 10:51:00.539843'),
    @('ATL3', 'atlastin GTPase 3', '3', '2021-10-05 10:51:00.539846'),
    @('ATP1A1', 'ATPase Na+/K+ transporting subunit alpha 1', '3', '2021-10-05 10:51:00.539848'),
    @('ATP7A', 'ATPase copper transporting alpha', '3', '2021-10-05 10:51:00.539851'),
    @('BICD2', 'BICD cargo adaptor 2', '3', '2021-10-05 10:51:00.539854'),
    @('BSCL2', 'BSCL2, seipin lipid droplet biogenesis associated', '3', '2021-10-05 10:51:00.539856'),
    @('CHCHD10', 'coiled-coil-helix-coiled-coil-helix domain containing 10', '3', '2021-10-05 10:51:00.539859'),
    @('COX20', 'COX20, cytochrome c oxidase assembly factor', '3', '2021-10-05 10:51:00.539862'),
    @('COX6A1', 'cytochrome c oxidase subunit 6A1', '3', '2021-10-05 10:51:00.539864'),
    @('DCTN1', 'dynactin subunit 1', '3', '2021-10-05 10:51:00.539867'),
    @('DNAJB2', 'DnaJ heat shock protein family (Hsp40) member B2', '3', '2021-10-05 10:51:00.539869'),
    @('DNM2', 'dynamin 2', '3', '2021-10-05 10:51:00.539872'),
    @('DNMT1', 'DNA methyltransferase 1', '3', '2021-10-05 10:51:00.539874'),
    @('DRP2', 'dystrophin related protein 2', '3', '2021-10-05 10:51:00.539877'),
    @('DST', 'dystonin', '3', '2021-10-05 10:51:00.539879'),
    @('DYNC1H1', 'dynein cytoplasmic 1 heavy chain 1', '3', '2021-10-05 10:51:00.539882'),
    @('EGR2', 'early growth response 2', '3', '2021-10-05 10:51:00.539884'),
    @('ELP1', 'elongator complex protein 1', '3', '2021-10-05 10:51:00.539887'),
    @('FBLN5', 'fibulin 5', '3', '2021-10-05 10:51:00.539890'),
    @('FGD4', 'FYVE, RhoGEF and PH domain containing 4', '3', '2021-10-05 10:51:00.539892'),
    @('FIG4', 'FIG4 phosphoinositide 5-phosphatase', '3', '2021-10-05 10:51:00.539894'),
    @('GARS', 'glycyl-tRNA synthetase', '3', '2021-10-05 10:51:00.539897'),
    @('GBF1', 'golgi brefeldin A resistant guanine nucleotide exchange factor 1', '3', '2021-10-05 10:51:00.539899'),
    @('GDAP1', 'ganglioside induced differentiation associated protein 1', '3', '2021-10-05 10:51:00.539902'),
    @('GJB1', 'gap junction protein beta 1', '3', '2021-10-05 10:51:00.539905'),
    @('GNB4', 'G protein subunit beta 4', '3', '2021-10-05 10:51:00.539908'),
    @('HARS', 'histidyl-tRNA synthetase', '3', '2021-10-05 10:51:00.539910'),
    @('HINT1', 'histidine triad nucleotide binding protein 1', '3', '2021-10-05 10:51:00.539913'),
    @('HK1', 'hexokinase 1', '3', '2021-10-05 10:51:00.539915'),
    @('HSPB1', 'heat shock protein family B (small) member 1', '3', '2021-10-05 10:51:00.539918'),
    @('HSPB8', 'heat shock protein family B (small) member 8', '3', '2021-10-05 10:51:00.539920'),
    @('IGHMBP2', 'immunoglobulin mu binding protein 2', '3', '2021-10-05 10:51:00.539923'),
    @('INF2', 'inverted formin, FH2 and WH2 domain containing', '3', '2021-10-05 10:51:00.539926'),
    @('JAG1', 'jagged 1', '3', '2021-10-05 10:51:00.539928'),
    @('KIF1A', 'kinesin family member 1A', '3', '2021-10-05 10:51:00.539931'),
    @('KIF5A', 'kinesin family member 5A', '3', '2021-10-05 10:51:00.539933'),
    @('LITAF', 'lipopolysaccharide induced TNF factor', '3', '2021-10-05 10:51:00.539936'),
    @('LRSAM1', 'leucine rich repeat and sterile alpha motif containing 1', '3', '2021-10-05 10:51:00.539938'),
    @('MFN2', 'mitofusin 2', '3', '2021-10-05 10:51:00.539941'),
    @('MME', 'membrane metalloendopeptidase', '3', '2021-10-05 10:51:00.539944'),
    @('MORC2', 'MORC family CW-type zinc finger 2', '3', '2021-10-05 10:51:00.539946'),
    @('MPV17', 'MPV17, mitochondrial inner membrane protein', '3', '2021-10-05 10:51:00.539949'),
    @('MPZ', 'myelin protein zero', '3', '2021-10-05 10:51:00.539951'),
    @('MTMR2', 'myotubularin related protein 2', '3', '2021-10-05 10:51:00.539954'),
    @('NDRG1', 'N-myc downstream regulated 1', '3', '2021-10-05 10:51:00.539956'),
    @('NEFH', 'neurofilament heavy', '3', '2021-10-05 10:51:00.539959'),
    @('NEFL', 'neurofilament light', '3', '2021-10-05 10:51:00.539961'),
    @('NGF', 'nerve growth factor', '3', '2021-10-05 10:51:00.539964'),
    @('PDK3', 'pyruvate dehydrogenase kinase 3', '3', '2021-10-05 10:51:00.539966'),
    @('PLEKHG5', 'pleckstrin homology and RhoGEF domain containing G5', '3', '2021-10-05 10:51:00.539969'),
    @('PMP2', 'peripheral myelin protein 2', '3', '2021-10-05 10:51:00.539972'),
    @('PMP22', 'peripheral myelin protein 22', '3', '2021-10-05 10:51:00.539974'),
    @('PRDM12', 'PR/SET domain 12', '3', '2021-10-05 10:51:00.539977'),
    @('PRPS1', 'phosphoribosyl pyrophosphate synthetase 1', '3', '2021-10-05 10:51:00.539980'),
    @('PRX', 'periaxin', '3', '2021-10-05 10:51:00.539982'),
    @('RAB7A', 'RAB7A, member RAS oncogene family', '3', '2021-10-05 10:51:00.539985'),
    @('REEP1', 'receptor accessory protein 1', '3', '2021-10-05 10:51:00.539987'),
    @('RETREG1', 'reticulophagy regulator 1', '3', '2021-10-05 10:51:00.539990'),
    @('SBF1', 'SET binding factor 1', '3', '2021-10-05 10:51:00.539992'),
    @('SBF2', 'SET binding factor 2', '3', '2021-10-05 10:51:00.539995'),
    @('SCN10A', 'sodium voltage-gated channel alpha subunit 10', '3', '2021-10-05 10:51:00.539997'),
    @('SCN11A', 'sodium voltage-gated channel alpha subunit 11', '3', '2021-10-05 10:51:00.540000'),
    @('SCN9A', 'sodium voltage-gated channel alpha subunit 9', '3', '2021-10-05 10:51:00.540002'),
    @('SEPT9', 'septin 9', '3', '2021-10-05 10:51:00.540005'),
    @('SH3TC2', 'SH3 domain and tetratricopeptide repeats 2', '3', '2021-10-05 10:51:00.540008'),
    @('SIGMAR1', 'sigma non-opioid intracellular receptor 1', '3', '2021-10-05 10:51:00.540011'),
    @('SLC5A7', 'solute carrier family 5 member 7', '3', '2021-10-05 10:51:00.540014'),
    @('SMN1', 'survival of motor neuron 1, telomeric', '3', '2021-10-05 10:51:00.540016'),
    @('SORD', 'sorbitol dehydrogenase', '3', '2021-10-05 10:51:00.540019'),
    @('SPG11', 'SPG11, spatacsin vesicle trafficking associated', '3', '2021-10-05 10:51:00.540021'),
    @('SPTAN1', 'spectrin alpha, non-erythrocytic 1', '3', '2021-10-05 10:51:00.540024'),
    @('SPTLC1', 'serine palmitoyltransferase long chain base subunit 1', '3', '2021-10-05 10:51:00.540026'),
    @('SPTLC2', 'serine palmitoyltransferase long chain base subunit 2', '3', '2021-10-05 10:51:00.540029'),
    @('SYT2', 'synaptotagmin 2', '3', '2021-10-05 10:51:00.540031'),
    @('TFG', 'TRK-fused gene', '3', '2021-10-05 10:51:00.540034'),
    @('TRIM2', 'tripartite motif containing 2', '3', '2021-10-05 10:51:00.540038'),
    @('TRPV4', 'transient receptor potential cation channel subfamily V member 4', '3', '2021-10-05 10:51:00.540041'),
    @('UBA1', 'ubiquitin like modifier activating enzyme 1', '3', '2021-10-05 10:51:00.540044'),
    @('VAPB', 'VAMP associated protein B and C', '3', '2021-10-05 10:51:00.540046'),
    @('VCP', 'valosin containing protein', '3', '2021-10-05 10:51:00.540049'),
    @('VRK1', 'vaccinia related kinase 1', '3', '2021-10-05 10:51:00.540051'),
    @('VWA1', 'von Willebrand factor A domain containing 1', '3', '2021-10-05 10:51:00.540054'),
    @('WARS', 'tryptophanyl-tRNA synthetase', '3', '2021-10-05 10:51:00.540056'),
    @('WNK1', 'WNK lysine deficient protein kinase 1', '3', '2021-10-05 10:51:00.540059'),
    @('YARS', 'tyrosyl-tRNA synthetase', '3', '2021-10-05 10:51:00.540062'),
    @('ARHGEF10', 'Rho guanine nucleotide exchange factor 10', '2', '2021-10-05 10:51:00.540064'),
    @('C1orf194', 'chromosome 1 open reading frame 194', '2', '2021-10-05 10:51:00.540067'),
    @('CADM3', 'cell adhesion molecule 3', '2', '2021-10-05 10:51:00.540070'),
    @('DGAT2', 'diacylglycerol O-acyltransferase 2', '2', '2021-10-05 10:51:00.540072'),
    @('DHTKD1', 'dehydrogenase E1 and transketolase domain containing 1', '2', '2021-10-05 10:51:00.540075'),
    @('FBXO38', 'F-box protein 38', '2', '2021-10-05 10:51:00.540077'),
    @('ITPR3', 'inositol 1,4,5-trisphosphate receptor type 3', '2', '2021-10-05 10:51:00.540081'),
    @('KIF1B', 'kinesin family member 1B', '2', '2021-10-05 10:51:00.540084'),
    @('LMNA', 'lamin A/C', '2', '2021-10-05 10:51:00.540087'),
    @('MARS', 'methionyl-tRNA synthetase', '2', '2021-10-05 10:51:00.540089'),
    @('NAGLU', 'N-acetyl-alpha-glucosaminidase', '2', '2021-10-05 10:51:00.540092'),
    @('RBM7', 'RNA binding motif protein 7', '2', '2021-10-05 10:51:00.540095'),
    @('SCO2', 'SCO2, cytochrome c oxidase assembly protein', '2', '2021-10-05 10:51:00.540097'),
    @('UBA5', 'ubiquitin like modifier activating enzyme 5', '2', '2021-10-05 10:51:00.540100'),
    @('HSPB3', 'heat shock protein family B (small) member 3', '1', '2021-10-05 10:51:00.540103'),
    @('IQGAP3', 'IQ motif containing GTPase activating protein 3', '1', '2021-10-05 10:51:00.540105'),
    @('KLHL13', 'kelch like family member 13', '1', '2021-10-05 10:51:00.540108'),
    @('LAS1L', 'LAS1 like, ribosome biogenesis factor', '1', '2021-10-05 10:51:00.540111'),
    @('MED25', 'mediator complex subunit 25', '1', '2021-10-05 10:51:00.540113'),
    @('SH3BP4', 'SH3 domain binding protein 4', '1', '2021-10-05 10:51:00.540116'),
    @('TRPA1', 'transient receptor potential cation channel subfamily A member 1', '1', '2021-10-05 10:51:00.540119')
)

$panelName = "Hereditary Neuropathy_CMT - isolated"

# New header cell F1 ("time_taken"), formatted like the other header cells
# (bold, centered, bordered) by copying the formatting from the existing E1 header.
$ws.Range("E1").Copy($ws.Range("F1")) | Out-Null
$ws.Range("F1").Value = "time_taken"

# Row 108 is entirely new (the sheet used to stop at row 107); give its index
# cell (A108) the same formatting as the other index cells in column A.
$ws.Range("A107").Copy($ws.Range("A108")) | Out-Null

# Columns D (geneConfidence) and F (time_taken) hold text that looks numeric /
# date-like ("3", "2021-10-05 10:51:00.539831", ...). Force a text format so
# Excel does not silently convert them to numbers / dates.
$ws.Range("D2:D108").NumberFormat = "@"
$ws.Range("F2:F108").NumberFormat = "@"

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]

    $ws.Cells.Item($r, 1).Value = $i
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $panelName
    $ws.Cells.Item($r, 6).Value = $row[3]
}

Write-Host "Rows written:" $data.Length
